$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# Row 3: I3 becomes a plain SUM instead of continuing the shared
# formula from B3:H3
# ---------------------------------------------------------------
$ws.Range("I3").Formula = '=SUM(B3:H3)'

# ---------------------------------------------------------------
# Row 19: fill in F/G/H + total column
# ---------------------------------------------------------------
$ws.Range("F19").Value = 5
$ws.Range("G19").Value = 5
$ws.Range("H19").Value = 5
$ws.Range("I19").Formula = '=SUM(B19:H19)'

# ---------------------------------------------------------------
# Row 20: new data for C..H + average column
# ---------------------------------------------------------------
$ws.Range("C20").Value = 0.04
$ws.Range("D20").Value = 0
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 0
$ws.Range("G20").Value = 0
$ws.Range("H20").Value = 0.01
$ws.Range("I20").Formula = '=AVERAGE(B20:H20)'

# ---------------------------------------------------------------
# Row 21: fill in F/G/H + total column
# ---------------------------------------------------------------
$ws.Range("F21").Value = 3
$ws.Range("G21").Value = 3
$ws.Range("H21").Value = 3
$ws.Range("I21").Formula = '=SUM(B21:H21)'

# ---------------------------------------------------------------
# Row 22: new data for D..H + average column
# ---------------------------------------------------------------
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 0.065
$ws.Range("F22").Value = 0
$ws.Range("G22").Value = 0
$ws.Range("H22").Value = 0
$ws.Range("I22").Formula = '=AVERAGE(B22:H22)'

# ---------------------------------------------------------------
# Row 23: extend the shared formula across F:H + total column
# ---------------------------------------------------------------
$ws.Range("F23:H23").Formula = '=(F19+F21)*$L$20'
$ws.Range("I23").Formula = '=SUM(B23:H23)'

# ---------------------------------------------------------------
# Row 24: apply the (1-row20) discount factor + total column
# ---------------------------------------------------------------
$ws.Range("C24").Formula = '=C2*$L$36*(1-C20)'
$ws.Range("D24").Formula = '=D2*$L$36*(1-D20)'
$ws.Range("E24").Formula = '=E2*$L$36*(1-E20)'
$ws.Range("F24").Formula = '=F2*$L$36*(1-F20)'
$ws.Range("G24").Formula = '=G2*$L$36*(1-G20)'
$ws.Range("H24").Formula = '=H2*$L$36*(1-H20)'
$ws.Range("I24").Formula = '=SUM(B24:H24)'

# ---------------------------------------------------------------
# Row 25: total column
# ---------------------------------------------------------------
$ws.Range("I25").Formula = '=SUM(B25:H25)'

# ---------------------------------------------------------------
# Row 26: recompute (values shift because row 22 changed) + total
# ---------------------------------------------------------------
$ws.Range("B26").Formula = '=(1-B22)*B2+B25'
$ws.Range("C26").Formula = '=(1-C22)*C2+C25'
$ws.Range("D26").Formula = '=(1-D22)*D2+D25'
$ws.Range("E26").Formula = '=(1-E22)*E2+E25'
$ws.Range("F26").Formula = '=(1-F22)*F2+F25'
$ws.Range("G26").Formula = '=(1-G22)*G2+G25'
$ws.Range("H26").Formula = '=(1-H22)*H2+H25'
$ws.Range("I26").Formula = '=SUM(B26:H26)'

# ---------------------------------------------------------------
# Row 27: drop the stray number format on B/F/H, refresh values,
# add total column
# ---------------------------------------------------------------
$ws.Range("B27").ClearFormats()
$ws.Range("F27").ClearFormats()
$ws.Range("H27").ClearFormats()
$ws.Range("B27").Formula = '=B23+(B24+B26)*$L$2*$N$2'
$ws.Range("C27").Formula = '=C23+(C24+C26)*$L$2*$N$2'
$ws.Range("D27").Formula = '=D23+(D24+D26)*$L$2*$N$2'
$ws.Range("E27").Formula = '=E23+(E24+E26)*$L$2*$N$2'
$ws.Range("F27").Formula = '=F23+(F24+F26)*$L$2*$N$2'
$ws.Range("G27").Formula = '=G23+(G24+G26)*$L$2*$N$2'
$ws.Range("H27").Formula = '=H23+(H24+H26)*$L$2*$N$2'
$ws.Range("I27").Formula = '=SUM(B27:H27)'

# ---------------------------------------------------------------
# Row 28: refresh values, add total column
# ---------------------------------------------------------------
$ws.Range("B28").Formula = '=B27*$L$25'
$ws.Range("C28").Formula = '=C27*$L$25'
$ws.Range("D28").Formula = '=D27*$L$25'
$ws.Range("E28").Formula = '=E27*$L$25'
$ws.Range("F28").Formula = '=F27*$L$25'
$ws.Range("G28").Formula = '=G27*$L$25'
$ws.Range("H28").Formula = '=H27*$L$25'
$ws.Range("I28").Formula = '=SUM(B28:H28)'

# ---------------------------------------------------------------
# Selection moves to F32 (matches the saved cursor position)
# ---------------------------------------------------------------
$ws.Range("F32").Select()
